$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot the original values of rows 59..86 (A:R) before we touch anything.
# new row N (60..86) will receive the values currently sitting in row (N-1) (59..85).
# new row 87 (inserted) will receive the values currently sitting in row 86.
$snapshot = @{}
for ($r = 59; $r -le 86; $r++) {
    $snapshot[$r] = $ws.Range("A$r`:R$r").Value2
}

# --- Step 2: insert a new row at position 87 (old row 87 shifts down to become row 88).
$ws.Rows.Item(87).Insert()

# --- Step 3: fill the newly inserted row 87 with the snapshot of the old row 86.
$ws.Range("A87:R87").Value = $snapshot[86]

# --- Step 4: shift rows 60..86 down from the snapshot of the row directly above them.
for ($r = 86; $r -ge 60; $r--) {
    $ws.Range("A$r`:R$r").Value = $snapshot[$r - 1]
}

# --- Step 5: row 59 gets brand-new data (new market day), reusing the unchanged columns.
$ws.Range("I59").Value = "Primera"
$ws.Range("D59").Value = 45029
$ws.Range("J59").Value = 40
$ws.Range("L59").Value = 20000
$ws.Range("M59").Value = 19000
$ws.Range("P59").Value = 1056
